$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A handful of Price cells now hold plain decimal-looking text (e.g. "61.67").
# The sheet stores every Price/Volume cell as text, so force those specific
# cells to the Text number format first -- otherwise Excel would silently
# re-interpret the assigned string as a numeric value.
$textCells = @("D17", "D18", "D23", "D25", "D26", "D27", "D34", "D36", "D37", "D39", "D43", "D45", "D48", "D49")
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = "26.911.59"
$ws.Range("E2").Value = "  +0.69%  "
$ws.Range("D3").Value = "1.552.67"
$ws.Range("E3").Value = "  +1.21%  "
$ws.Range("E4").Value = "  +0.63%  "
$ws.Range("E5").Value = "  +0.49%  "
$ws.Range("E7").Value = "  +0.58%  "
$ws.Range("E8").Value = "  +1.60%  "
$ws.Range("E9").Value = "  +1.39%  "
$ws.Range("E10").Value = "  +1.02%  "
$ws.Range("E11").Value = "  +0.60%  "
$ws.Range("D12").Value = "1.773.97"
$ws.Range("E12").Value = "  +1.27%  "
$ws.Range("D13").Value = "1.527.03"
$ws.Range("E13").Value = "  -0.46%  "
$ws.Range("E14").Value = "  +1.50%  "
$ws.Range("E15").Value = "  +1.88%  "
$ws.Range("D16").Value = "26.915.54"
$ws.Range("E16").Value = "  +0.74%  "
$ws.Range("D17").Value = "61.67"
$ws.Range("E17").Value = "  +1.03%  "
$ws.Range("D18").Value = "217.16"
$ws.Range("E18").Value = "  +2.18%  "
$ws.Range("E19").Value = "  +0.99%  "
$ws.Range("E20").Value = "  -0.11%  "
$ws.Range("E21").Value = "  +0.57%  "
$ws.Range("E22").Value = "  +0.89%  "
$ws.Range("D23").Value = "9.23"
$ws.Range("E23").Value = "  +1.63%  "
$ws.Range("E24").Value = "  +0.51%  "
$ws.Range("D25").Value = "153.86"
$ws.Range("E25").Value = "  +1.34%  "
$ws.Range("D26").Value = "6.58"
$ws.Range("D27").Value = "14.86"
$ws.Range("E27").Value = "  +0.42%  "
$ws.Range("E28").Value = "  +0.65%  "
$ws.Range("E29").Value = "  +1.03%  "
$ws.Range("E30").Value = "  +2.82%  "
$ws.Range("E31").Value = "  +0.08%  "
$ws.Range("D33").Value = "1.420.19"
$ws.Range("E33").Value = "  +4.19%  "
$ws.Range("D34").Value = "3.01"
$ws.Range("E34").Value = "  +2.96%  "
$ws.Range("E35").Value = "  +3.88%  "
$ws.Range("D36").Value = "0.958"
$ws.Range("E36").Value = "  +1.44%  "
$ws.Range("D37").Value = "2.30"
$ws.Range("E37").Value = "  +1.04%  "
$ws.Range("E38").Value = "  +0.89%  "
$ws.Range("D39").Value = "0.521"
$ws.Range("E39").Value = "  +0.04%  "
$ws.Range("E40").Value = "  +1.32%  "
$ws.Range("E41").Value = "  +0.61%  "
$ws.Range("E42").Value = "  -0.24%  "
$ws.Range("D43").Value = "0.987"
$ws.Range("E43").Value = "  -0.89%  "
$ws.Range("D45").Value = "63.64"
$ws.Range("E45").Value = "  +1.78%  "
$ws.Range("E46").Value = "  -0.70%  "
$ws.Range("D47").Value = "1.688.58"
$ws.Range("E47").Value = "  +1.45%  "
$ws.Range("D48").Value = "86.23"
$ws.Range("E48").Value = "  +0.97%  "
$ws.Range("D49").Value = "0.0523"
$ws.Range("E49").Value = "  +4.19%  "
$ws.Range("D50").Value = "0.0₆0102"
$ws.Range("E50").Value = "  +4.29%  "
$ws.Range("E51").Value = "  +1.50%  "
